$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"1.131426898102946E-15"
$ws.Range("E2").Value = [double]"1.131426898102946E-15"

$ws.Range("D3").Value = [double]"2.984313054960899E-12"
$ws.Range("E3").Value = [double]"2.984313054960899E-12"

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0

$ws.Range("D5").Value = [double]"5.605916941340763E-08"
$ws.Range("E5").Value = 0.9999999439408306

$ws.Range("D6").Value = [double]"6.998663586082739E-15"
$ws.Range("E6").Value = 0.999999999999993

$ws.Range("D7").Value = [double]"2.497628555398735E-10"
$ws.Range("E7").Value = 0.9999999997502371

$ws.Range("D8").Value = [double]"1.307942752328771E-10"
$ws.Range("E8").Value = 0.9999999998692057
$ws.Range("F8").Value = 13.45111751556396
